$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N2").Value = 0.1265680465041342
$ws.Range("N3").Value = 0.02271239078883044
$ws.Range("N4").Value = 0.01279536440308497
$ws.Range("N5").Value = 0.001
$ws.Range("N6").Value = 0.02873139331440719
$ws.Range("N7").Value = 0.04739794823921178
$ws.Range("N8").Value = 0.15
$ws.Range("N9").Value = 0.001000000000000011
$ws.Range("N10").Value = 0.15
$ws.Range("N11").Value = 0.001
$ws.Range("N12").Value = 0.001
$ws.Range("N13").Value = 0.15
$ws.Range("N14").Value = 0.001
$ws.Range("N15").Value = 0.08745057021659455
$ws.Range("N16").Value = 0.001000000000000005
$ws.Range("N17").Value = 0.001000000000000013
$ws.Range("N18").Value = 0.15
$ws.Range("N19").Value = 0.001000000000000009
$ws.Range("N20").Value = 0.06534428653373688
$ws.Range("N21").Value = 0.001000000000000003
$ws.Range("B22").Value = -0.03755441501883384
$ws.Range("C22").Value = -0.0172326475632052
$ws.Range("D22").Value = 0.02873724367175603
$ws.Range("E22").Value = 0.005595375681225705
$ws.Range("F22").Value = 0.07494363704725722
$ws.Range("G22").Value = -0.01347926687622574
$ws.Range("H22").Value = -0.01526264754209035
$ws.Range("I22").Value = 0.03589284412162125
$ws.Range("J22").Value = -0.002528546400419587
$ws.Range("K22").Value = 0.01595208085631053
$ws.Range("L22").Value = 0.02418446431986731
$ws.Range("M22").Value = -0.008285889085776309
$ws.Range("B23").Value = 0.9631420069044215
$ws.Range("C23").Value = 0.9829149852566068
$ws.Range("D23").Value = 1.029154142181564
$ws.Range("E23").Value = 1.005611059033504
$ws.Range("F23").Value = 1.077823399863242
$ws.Range("G23").Value = 0.9866111716372192
$ws.Range("H23").Value = 0.9848532363488162
$ws.Range("I23").Value = 1.036544768674543
$ws.Range("J23").Value = 0.9974746476803357
$ws.Range("K23").Value = 1.016079994556233
$ws.Range("L23").Value = 1.024479280335336
$ws.Range("M23").Value = 0.9917483442766625
$ws.Range("N23").Value = 1.129412494234728
